$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Se cambio la edad de ana
$ws.Range("B2").Value = 24

# Se cambio la provincia de luis
$ws.Range("C3").Value = "Guanacaste"

# Update selection to match target state
$ws.Range("D12").Select()
